$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add three new parameter rows (94-96) to the table, following the same
# Name / Item(code) / Value / index layout used by the rest of the sheet:
#   94: Passengers capcity public system | capacity-public   | 136  | 26
#   95: Regular-speed-public             | reg-speed-pub     | 0.15 | 27
#   96: Default waiting time in the public system
#                                         | default-wait-time | 6.6  | 28
# ---------------------------------------------------------------------

# Column D on this sheet is always styled the same way (centered,
# default font/format). Grab that look from the row above and stamp it
# down over the three new D cells, then fill in the running index.
$ws.Range("D93").Copy()
$ws.Range("D94:D96").PasteSpecial(-4122)
$ws.Range("D94").Value = 26
$ws.Range("D95").Value = 27
$ws.Range("D96").Value = 28

# Row 94's C cell is a plain number, styled like the other numeric
# entries in column C (e.g. C84).
$ws.Range("C84").Copy()
$ws.Range("C94").PasteSpecial(-4122)
$ws.Range("C94").Value = 136

# Rows 95/96's C cells hold numbers typed as text ("0.15" / "6.6"),
# but keep the same plain style used elsewhere in column D/C94 (no
# special text format). Build the text via a throwaway formula cell so
# the value comes back as a real string, then drop just that value
# (not its format) onto the already-styled destination cells.
$ws.Range("D94").Copy()
$ws.Range("C95:C96").PasteSpecial(-4122)

$ws.Range("ZZ1").Formula = "=""0.15"""
$ws.Range("ZZ1").Copy()
$ws.Range("C95").PasteSpecial(-4163)

$ws.Range("ZZ1").Formula = "=""6.6"""
$ws.Range("ZZ1").Copy()
$ws.Range("C96").PasteSpecial(-4163)

$ws.Range("ZZ1").Clear()

# Column A / B text labels + codes for the new rows.
$ws.Range("A94").Value = "Passengers capcity public system"
$ws.Range("B94").Value = "capacity-public"

$ws.Range("A95").Value = "Regular-speed-public"
$ws.Range("B95").Value = "reg-speed-pub"

$ws.Range("A96").Value = "Default waiting time in the public system"
$ws.Range("B96").Value = "default-wait-time"

# Reflect where the user ended up after entering the new data.
$ws.Range("C94:C96").Select()
